$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $value) {
    # Force the cell to be written as a text/string value (matching the
    # workbook's existing inline-string cells) instead of letting Excel
    # auto-coerce numeric-looking text into a real number.
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

# --- Rows 23/24 swapped places (Dai <-> Polygon) and got new data ---
Set-TextValue 23 2 "Polygon"
Set-TextValue 23 3 "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue 23 4 "0.541"
Set-TextValue 23 5 "  +1.83%  "

Set-TextValue 24 2 "Dai"
Set-TextValue 24 3 "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue 24 4 "0.999"
Set-TextValue 24 5 "  -0.03%  "

# --- Price (column D) / Volume(1h) (column E) refresh for the remaining rows ---
$data = @(
    @{Row=2; D="67.324.61"; E="  +0.52%  "},
    @{Row=3; D="3.522.87"; E="  +0.57%  "},
    @{Row=4; E="  -0.07%  "},
    @{Row=5; D="596.56"; E="  +0.51%  "},
    @{Row=6; D="173.40"; E="  +2.53%  "},
    @{Row=7; D="0.999"; E="  -0.02%  "},
    @{Row=8; D="0.594"; E="  +2.94%  "},
    @{Row=9; D="0.135"; E="  +7.84%  "},
    @{Row=10; E="  +0.39%  "},
    @{Row=11; D="0.437"; E="  -0.69%  "},
    @{Row=12; D="4.130.44"; E="  +0.31%  "},
    @{Row=13; D="0.135"; E="  -0.11%  "},
    @{Row=14; D="28.73"; E="  +1.86%  "},
    @{Row=15; D="0.0000182"; E="  +1.60%  "},
    @{Row=16; D="67.293.73"; E="  +0.55%  "},
    @{Row=17; D="3.510.84"; E="  -0.52%  "},
    @{Row=18; D="6.34"; E="  +0.28%  "},
    @{Row=19; D="14.29"; E="  +1.50%  "},
    @{Row=20; D="398.29"; E="  +0.65%  "},
    @{Row=21; D="8.01"; E="  +1.11%  "},
    @{Row=22; D="73.47"; E="  -0.24%  "},
    @{Row=25; D="0.0000123"; E="  -3.80%  "},
    @{Row=26; D="10.23"; E="  +1.42%  "},
    @{Row=28; D="0.998"; E="  -0.22%  "},
    @{Row=29; D="6.29"; E="  -1.64%  "},
    @{Row=30; E="  -0.86%  "},
    @{Row=31; D="2.09"; E="  +1.12%  "},
    @{Row=32; D="24.14"; E="  +2.42%  "},
    @{Row=33; D="7.40"},
    @{Row=34; D="1.63"; E="  +2.96%  "},
    @{Row=35; D="164.04"; E="  +1.21%  "},
    @{Row=36; D="0.896"; E="  -1.10%  "},
    @{Row=37; D="1.92"; E="  -1.23%  "},
    @{Row=38; D="6.92"; E="  +3.39%  "},
    @{Row=39; D="4.73"; E="  +1.18%  "},
    @{Row=40; D="0.0746"; E="  -0.87%  "},
    @{Row=41; D="27.44"; E="  +2.98%  "},
    @{Row=42; D="26.44"; E="  +0.27%  "},
    @{Row=43; D="2.827.91"; E="  -0.12%  "},
    @{Row=44; D="2.63"; E="  +4.30%  "},
    @{Row=45; D="42.91"; E="  -1.44%  "},
    @{Row=46; E="  -2.78%  "},
    @{Row=47; D="341.32"; E="  -2.37%  "},
    @{Row=48; E="  +1.09%  "},
    @{Row=49; D="33.76"; E="  +1.07%  "},
    @{Row=50; D="6.55"; E="  +0.21%  "},
    @{Row=51; D="0.853"; E="  +0.16%  "}
)

foreach ($item in $data) {
    if ($item.ContainsKey("D")) { Set-TextValue $item.Row 4 $item.D }
    if ($item.ContainsKey("E")) { Set-TextValue $item.Row 5 $item.E }
}
